$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 168
$ws1.Range("F10").Value = 5446
$ws1.Range("F11").Value = 4865

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 78
$ws2.Range("G4").Value = 64

# Sheet "全部类型" (All types) - combined view
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 168
$ws4.Range("F10").Value = 5446
$ws4.Range("F11").Value = 4865
$ws4.Range("F17").Value = 78
$ws4.Range("G19").Value = 64
